$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 data (mirrors the existing pattern used by rows 3-6: columns A-G)
# A7 looks numeric ("1234567890") so force it to stay text, matching A3/A5,
# then reset the style so no numeric/quote-prefix formatting sticks around.
$ws.Cells.Item(7, 1).Value = "'1234567890"
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(7, 2).Value = "car"
$ws.Cells.Item(7, 3).Value = "red"
$ws.Cells.Item(7, 4).Value = "Toyota"
$ws.Cells.Item(7, 5).Value = "John Doe"
$ws.Cells.Item(7, 6).Value = 1768662257.194566
$ws.Cells.Item(7, 7).Value = "1234567890-4601"

# Row 8 data
$ws.Cells.Item(8, 1).Value = "test car"
$ws.Cells.Item(8, 2).Value = "car"
$ws.Cells.Item(8, 3).Value = "red"
$ws.Cells.Item(8, 4).Value = "Toyota"
$ws.Cells.Item(8, 5).Value = "John Doe"
$ws.Cells.Item(8, 6).Value = 1768662257.230299
$ws.Cells.Item(8, 7).Value = "test car-2789"
